$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.785.55"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "1.872.85"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7185"
$ws.Range("E5").Value = "  -2.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.79"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3145"
$ws.Range("E8").Value = "  -0.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07543"
$ws.Range("E9").Value = "  +5.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.58"
$ws.Range("E10").Value = "  -0.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08186"
$ws.Range("E11").Value = "  -2.09%  "
$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7439"
$ws.Range("E12").Value = "  -0.83%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.887.04"
$ws.Range("E13").Value = "  +0.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.335"
$ws.Range("E14").Value = "  -1.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.43"
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").Value = "29.738.87"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.019"
$ws.Range("E17").Value = "  -0.82%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "246.36"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007929"
$ws.Range("E19").Value = "  +1.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.46"
$ws.Range("E20").Value = "  -0.67%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9998"
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.113.30"
$ws.Range("E22").Value = "  -0.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.714"
$ws.Range("E24").Value = "  -3.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1501"
$ws.Range("E25").Value = "  -3.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.183"
$ws.Range("E26").Value = "  -0.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.76"
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.56"
$ws.Range("E28").Value = "  -0.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.003"
$ws.Range("E29").Value = "  -1.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.428"
$ws.Range("E30").Value = "  -5.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.552"
$ws.Range("E31").Value = "  -0.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.524"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.177"
$ws.Range("E33").Value = "  -2.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05446"
$ws.Range("E34").Value = "  +2.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.224"
$ws.Range("E35").Value = "  -1.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7372"
$ws.Range("E36").Value = "  -2.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9984"
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.703"
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01911"
$ws.Range("E39").Value = "  -2.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.730"
$ws.Range("E40").Value = "  -0.79%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4454"
$ws.Range("E41").Value = "  -1.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8914"
$ws.Range("E42").Value = "  +4.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.995"
$ws.Range("E43").Value = "  -0.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "71.61"
$ws.Range("E44").Value = "  -1.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.002"
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("D46").Value = "1.037.80"
$ws.Range("E46").Value = "  -6.80%  "
$ws.Range("E47").Value = "  +0.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.465"
$ws.Range("E48").Value = "  -2.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.799"
$ws.Range("E49").Value = "  -2.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.569"
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("D51").Value = "2.012.49"
$ws.Range("E51").Value = "  -0.65%  "
